$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Product_id"
$ws.Range("B1").Value = "Product_name"

$ws.Range("B1").Select()
